$d = $word.ActiveDocument

# --- Edit 1: first paragraph -----------------------------------------
# "This is a Microsoft word document." becomes
# "This is a Microsoft word document.  (This is a change – Version for main branch)"
# with the parenthetical annotation in red (FF0000), appended as three runs.
$firstPara = $d.Paragraphs.Item(1)
$firstRange = $firstPara.Range

# Insert two trailing spaces right after the existing sentence, before the
# paragraph mark at the end of the paragraph range.
$insPoint = $d.Range($firstRange.End - 1, $firstRange.End - 1)
$insPoint.InsertAfter("  ")

# Append the red annotation text in three runs (matches the source edit).
$enDash = [char]0x2013
$r1 = $d.Range($insPoint.End, $insPoint.End)
$r1.InsertAfter("(This is a change " + $enDash + " Ve")
$r1.Font.Color = 255

$r2 = $d.Range($r1.End, $r1.End)
$r2.InsertAfter("rsion for main branch")
$r2.Font.Color = 255

$r3 = $d.Range($r2.End, $r2.End)
$r3.InsertAfter(")")
$r3.Font.Color = 255

# --- Edit 2: drop the trailing "ank God almighty, we are free at last." --
# paragraph (the final paragraph in the document, styled NormalWeb).
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
if ($lastPara.Range.Text -like "*God almighty*") {
    $lastPara.Range.Delete()
}
